# Auto-generated Excel COM-interop script
# Implements: insert one duplicate row at the end of each of the last 12
# timestamp-groups (sizes 8,8,8,6,6,6,4,4,4,2,2,2) in sheets 1,2,4, then
# append 3 brand-new single-row groups (next day readings) at the tail.

$wb = $excel.ActiveWorkbook

# ---- Sheet index 1 ----
$ws = $wb.Worksheets.Item(1)

$payload = New-Object 'object[,]' 1,9
$payload[0,1] = "0x01,0x90"
$payload[0,2] = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$payload[0,3] = "0x01,0x90,"
$payload[0,4] = "0xd"
$payload[0,5] = [double]"400"
$payload[0,6] = [double]"5.686312626471138e+23"
$payload[0,7] = [double]"400"
$payload[0,8] = [double]"13"

# Insert one duplicate row after the end of each tail group (process
# bottom-up so earlier insert points stay valid row numbers).
$ws.Rows("530:530").Insert()
$ws.Cells.Item(530,1).Value2 = [double]"45726.73158645834"
$ws.Range("B530:I530").Value2 = $payload
$ws.Rows("528:528").Insert()
$ws.Cells.Item(528,1).Value2 = [double]"45726.73156329861"
$ws.Range("B528:I528").Value2 = $payload
$ws.Rows("526:526").Insert()
$ws.Cells.Item(526,1).Value2 = [double]"45726.73154126157"
$ws.Range("B526:I526").Value2 = $payload
$ws.Rows("524:524").Insert()
$ws.Cells.Item(524,1).Value2 = [double]"45726.23144357639"
$ws.Range("B524:I524").Value2 = $payload
$ws.Rows("520:520").Insert()
$ws.Cells.Item(520,1).Value2 = [double]"45726.23142038195"
$ws.Range("B520:I520").Value2 = $payload
$ws.Rows("516:516").Insert()
$ws.Cells.Item(516,1).Value2 = [double]"45726.23139893518"
$ws.Range("B516:I516").Value2 = $payload
$ws.Rows("512:512").Insert()
$ws.Cells.Item(512,1).Value2 = [double]"45725.73130123843"
$ws.Range("B512:I512").Value2 = $payload
$ws.Rows("506:506").Insert()
$ws.Cells.Item(506,1).Value2 = [double]"45725.73127832176"
$ws.Range("B506:I506").Value2 = $payload
$ws.Rows("500:500").Insert()
$ws.Cells.Item(500,1).Value2 = [double]"45725.73125641204"
$ws.Range("B500:I500").Value2 = $payload
$ws.Rows("494:494").Insert()
$ws.Cells.Item(494,1).Value2 = [double]"45725.23115855324"
$ws.Range("B494:I494").Value2 = $payload
$ws.Rows("486:486").Insert()
$ws.Cells.Item(486,1).Value2 = [double]"45725.23113506944"
$ws.Range("B486:I486").Value2 = $payload
$ws.Rows("478:478").Insert()
$ws.Cells.Item(478,1).Value2 = [double]"45725.23111297454"
$ws.Range("B478:I478").Value2 = $payload

# Append 3 brand-new rows (next day readings) at the tail.
$appendStart = $ws.UsedRange.Rows.Count
$r = $appendStart + 1
$ws.Cells.Item($r,1).Value2 = [double]"45727.23168409722"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload
$r = $appendStart + 2
$ws.Cells.Item($r,1).Value2 = [double]"45727.23170618056"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload
$r = $appendStart + 3
$ws.Cells.Item($r,1).Value2 = [double]"45727.23172934028"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload

# ---- Sheet index 2 ----
$ws = $wb.Worksheets.Item(2)

$payload = New-Object 'object[,]' 1,9
$payload[0,1] = "0x01,0x90"
$payload[0,2] = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$payload[0,3] = "0x01,0x90,"
$payload[0,4] = "0xe"
$payload[0,5] = [double]"400"
$payload[0,6] = [double]"5.686312626471138e+23"
$payload[0,7] = [double]"400"
$payload[0,8] = [double]"14"

# Insert one duplicate row after the end of each tail group (process
# bottom-up so earlier insert points stay valid row numbers).
$ws.Rows("544:544").Insert()
$ws.Cells.Item(544,1).Value2 = [double]"45726.57956549768"
$ws.Range("B544:I544").Value2 = $payload
$ws.Rows("542:542").Insert()
$ws.Cells.Item(542,1).Value2 = [double]"45726.57954204861"
$ws.Range("B542:I542").Value2 = $payload
$ws.Rows("540:540").Insert()
$ws.Cells.Item(540,1).Value2 = [double]"45726.57952"
$ws.Range("B540:I540").Value2 = $payload
$ws.Rows("538:538").Insert()
$ws.Cells.Item(538,1).Value2 = [double]"45726.07942256945"
$ws.Range("B538:I538").Value2 = $payload
$ws.Rows("534:534").Insert()
$ws.Cells.Item(534,1).Value2 = [double]"45726.07939922454"
$ws.Range("B534:I534").Value2 = $payload
$ws.Rows("530:530").Insert()
$ws.Cells.Item(530,1).Value2 = [double]"45726.07937777778"
$ws.Range("B530:I530").Value2 = $payload
$ws.Rows("526:526").Insert()
$ws.Cells.Item(526,1).Value2 = [double]"45725.57928042824"
$ws.Range("B526:I526").Value2 = $payload
$ws.Rows("520:520").Insert()
$ws.Cells.Item(520,1).Value2 = [double]"45725.57925716435"
$ws.Range("B520:I520").Value2 = $payload
$ws.Rows("514:514").Insert()
$ws.Cells.Item(514,1).Value2 = [double]"45725.57923533564"
$ws.Range("B514:I514").Value2 = $payload
$ws.Rows("508:508").Insert()
$ws.Cells.Item(508,1).Value2 = [double]"45725.07913833333"
$ws.Range("B508:I508").Value2 = $payload
$ws.Rows("500:500").Insert()
$ws.Cells.Item(500,1).Value2 = [double]"45725.07911518519"
$ws.Range("B500:I500").Value2 = $payload
$ws.Rows("492:492").Insert()
$ws.Cells.Item(492,1).Value2 = [double]"45725.07909302083"
$ws.Range("B492:I492").Value2 = $payload

# Append 3 brand-new rows (next day readings) at the tail.
$appendStart = $ws.UsedRange.Rows.Count
$r = $appendStart + 1
$ws.Cells.Item($r,1).Value2 = [double]"45727.07966211806"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload
$r = $appendStart + 2
$ws.Cells.Item($r,1).Value2 = [double]"45727.07968420139"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload
$r = $appendStart + 3
$ws.Cells.Item($r,1).Value2 = [double]"45727.07970724537"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload

# ---- Sheet index 4 ----
$ws = $wb.Worksheets.Item(4)

$payload = New-Object 'object[,]' 1,9
$payload[0,1] = "0x01,0x90"
$payload[0,2] = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$payload[0,3] = "0x01,0x90,"
$payload[0,4] = "0x3"
$payload[0,5] = [double]"400"
$payload[0,6] = [double]"5.686312626471138e+23"
$payload[0,7] = [double]"400"
$payload[0,8] = [double]"3"

# Insert one duplicate row after the end of each tail group (process
# bottom-up so earlier insert points stay valid row numbers).
$ws.Rows("548:548").Insert()
$ws.Cells.Item(548,1).Value2 = [double]"45726.72967174769"
$ws.Range("B548:I548").Value2 = $payload
$ws.Rows("546:546").Insert()
$ws.Cells.Item(546,1).Value2 = [double]"45726.72964848379"
$ws.Range("B546:I546").Value2 = $payload
$ws.Rows("544:544").Insert()
$ws.Cells.Item(544,1).Value2 = [double]"45726.72962616898"
$ws.Range("B544:I544").Value2 = $payload
$ws.Rows("542:542").Insert()
$ws.Cells.Item(542,1).Value2 = [double]"45726.22953034722"
$ws.Range("B542:I542").Value2 = $payload
$ws.Rows("538:538").Insert()
$ws.Cells.Item(538,1).Value2 = [double]"45726.22950657408"
$ws.Range("B538:I538").Value2 = $payload
$ws.Rows("534:534").Insert()
$ws.Cells.Item(534,1).Value2 = [double]"45726.2294844213"
$ws.Range("B534:I534").Value2 = $payload
$ws.Rows("530:530").Insert()
$ws.Cells.Item(530,1).Value2 = [double]"45725.7293875"
$ws.Range("B530:I530").Value2 = $payload
$ws.Rows("524:524").Insert()
$ws.Cells.Item(524,1).Value2 = [double]"45725.72936453704"
$ws.Range("B524:I524").Value2 = $payload
$ws.Rows("518:518").Insert()
$ws.Cells.Item(518,1).Value2 = [double]"45725.7293421412"
$ws.Range("B518:I518").Value2 = $payload
$ws.Rows("512:512").Insert()
$ws.Cells.Item(512,1).Value2 = [double]"45725.22924497685"
$ws.Range("B512:I512").Value2 = $payload
$ws.Rows("504:504").Insert()
$ws.Cells.Item(504,1).Value2 = [double]"45725.22922125"
$ws.Range("B504:I504").Value2 = $payload
$ws.Rows("496:496").Insert()
$ws.Cells.Item(496,1).Value2 = [double]"45725.22919952546"
$ws.Range("B496:I496").Value2 = $payload

# Append 3 brand-new rows (next day readings) at the tail.
$appendStart = $ws.UsedRange.Rows.Count
$r = $appendStart + 1
$ws.Cells.Item($r,1).Value2 = [double]"45727.22976834491"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload
$r = $appendStart + 2
$ws.Cells.Item($r,1).Value2 = [double]"45727.22979072916"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload
$r = $appendStart + 3
$ws.Cells.Item($r,1).Value2 = [double]"45727.22981358796"
$ws.Cells.Item($r,1).Style = $ws.Cells.Item($r-1,1).Style
$ws.Range("B" + $r + ":I" + $r).Value2 = $payload

